# Insert a new weekly price record as row 303, pushing the existing
# rows 303:374 down to 304:375 (dimension grows from A1:R374 to A1:R375).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("303:303").Insert()

$ws.Cells.Item(303, 1).Value = 8
$ws.Cells.Item(303, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(303, 3).Value = "Coquimbo"
$ws.Cells.Item(303, 4).Value = 44785
$ws.Cells.Item(303, 5).Value = 4
$ws.Cells.Item(303, 6).Value = 100114013
$ws.Cells.Item(303, 7).Value = "Zanahoria"
$ws.Cells.Item(303, 8).Value = "Sin especificar"
$ws.Cells.Item(303, 9).Value = "Primera"
$ws.Cells.Item(303, 10).Value = 600
$ws.Cells.Item(303, 11).Value = 10000
$ws.Cells.Item(303, 12).Value = 11000
$ws.Cells.Item(303, 13).Value = 10500
$ws.Cells.Item(303, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(303, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(303, 16).Value = 525
$ws.Cells.Item(303, 17).Value = 20
$ws.Cells.Item(303, 18).Value = "Hortaliza"
